$d = $word.ActiveDocument

# Fix the spelling of "Nieuwenhuys" -> "Nieuwenhuijs" (author's correct surname)
$d.Content.Find.Execute("Nieuwenhuys", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Nieuwenhuijs", 2)
